# Apply the updated bilibili-scraped numbers/text from the gh-pages data refresh.
$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 941
$ws1.Range("F6").Value = 4901
$ws1.Range("F7").Value = 377
$ws1.Range("F8").Value = 556
$ws1.Range("F9").Value = 859
$ws1.Range("F12").Value = 14
$ws1.Range("D16").Value = "沈半路171号 T-Car杭州汽车文化主题公园"
$ws1.Range("F16").Value = 1582
$ws1.Range("I16").Value = "//i0.hdslb.com/bfs/openplatform/202402/7JXZatUK1707187527932.jpeg"
$ws1.Range("F18").Value = 676
$ws1.Range("F21").Value = 247
$ws1.Range("F27").Value = 1453
$ws1.Range("F28").Value = 132
$ws1.Range("F31").Value = 187
$ws1.Range("F36").Value = 251
$ws1.Range("F37").Value = 539
$ws1.Range("F38").Value = 70
$ws1.Range("G39").Value = 66

# --- Sheet "演出" (Show) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 129

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 941
$ws4.Range("F8").Value = 4901
$ws4.Range("F9").Value = 377
$ws4.Range("F10").Value = 556
$ws4.Range("F12").Value = 129
$ws4.Range("F13").Value = 859
$ws4.Range("F18").Value = 14
$ws4.Range("D23").Value = "沈半路171号 T-Car杭州汽车文化主题公园"
$ws4.Range("F23").Value = 1582
$ws4.Range("I23").Value = "//i0.hdslb.com/bfs/openplatform/202402/7JXZatUK1707187527932.jpeg"
$ws4.Range("F25").Value = 676
$ws4.Range("F28").Value = 247
$ws4.Range("F34").Value = 1453
$ws4.Range("F35").Value = 132
$ws4.Range("F38").Value = 187
$ws4.Range("F42").Value = 251
$ws4.Range("F43").Value = 539
$ws4.Range("F44").Value = 70
$ws4.Range("G45").Value = 66
